# Rouet et al., 2005 was measured by ELISA+Saturation in the source data;
# correct the recorded method to Radioligand across all three binding-partner
# sheets (VEGFR1, VEGFR2, NRP1).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VEGFA165_VEGFR1")
$ws2 = $wb.Worksheets.Item("VEGFA165_VEGFR2")
$ws3 = $wb.Worksheets.Item("VEGFA165_NRP1")

# --- the actual content fix: Rouet et al., 2005 method -> Radioligand ---
$ws1.Range("B11").Value = "Radioligand"
$ws2.Range("B11").Value = "Radioligand"
$ws3.Range("B9").Value = "Radioligand"

# The VEGFR2 sheet's last row (Soker et al., 1996) had picked up a stray
# fill-enabled style that isn't used anywhere else in the workbook; clear it
# so the row matches the plain bottom-border styling used on the other
# sheets' closing rows.
$ws2.Range("A14:C14").Interior.Pattern = -4142  # xlNone

# --- restore the editor's last on-screen selection state ---
$ws1.Range("B12").Select()
$ws3.Range("B10").Select()

$ws2.Activate()
$ws2.Range("B12").Select()
